$p = $ppt.ActivePresentation

# The Sprint #3 demo deck is being finalized: the trailing
# "Improvements" scratch slide (feature-request / UI-UX notes that
# didn't make the cut) is dropped from the final deck.
$targetIndex = -1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -eq "Improvements") {
                $targetIndex = $i
            }
        }
    }
}

if ($targetIndex -eq -1) {
    # Fallback: if the marker shape can't be found for some reason,
    # just drop the last slide (that's where it lives today).
    $targetIndex = $p.Slides.Count
}

$p.Slides.Item($targetIndex).Delete()
